# Insert a new weekly price record at row 40 (Fruta / Agrícola del Norte S.A. de
# Arica - Plátano), shifting all subsequent rows down by one, and carrying the
# former last row (133) down to the new last row (134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:133 down to 41:134, leaving a fresh (format-inherited) row 40.
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with the new observation.
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 44459
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100108
$ws.Cells.Item(40, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(40, 9).Value = 100108006
$ws.Cells.Item(40, 10).Value = "Plátano"
$ws.Cells.Item(40, 11).Value = "Sin especificar"
$ws.Cells.Item(40, 12).Value = "Pintón"
$ws.Cells.Item(40, 13).Value = 120
$ws.Cells.Item(40, 14).Value = 20000
$ws.Cells.Item(40, 15).Value = 21000
$ws.Cells.Item(40, 16).Value = 20500
$ws.Cells.Item(40, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(40, 18).Value = "Ecuador"
$ws.Cells.Item(40, 19).Value = 1025
$ws.Cells.Item(40, 20).Value = 20
